$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the header row (row 1) so that all the data rows shift up by one.
$ws.Rows.Item(1).Delete()

# Update the last-selected cell to reflect the new layout (G53 instead of H53).
$null = $ws.Range("G53").Select()
